$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (H1/I1), reusing the existing bold header style (s=1) ---
$ws.Range("A1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)
$ws.Range("H1").Value = "employeeid"
$ws.Range("I1").Value = "name"

# --- Replace row 2 data with the new single time-entry row ---
$ws.Range("B2").Value = "2024-05-27 00:00:00"
$ws.Range("C2").Value = "12:04:58"
$ws.Range("D2").Value = "0:00:00"
$ws.Range("E2").Value = "2024-05-27 12:05:15"
$ws.Range("F2").Value = "2024-05-27 12:05:15"
$ws.Range("G2").Value = "None"
$ws.Range("I2").Value = "alexa rodrig"

# A2 and H2 hold numeric-looking text ("12" / "3") that must stay text, not be
# coerced to numbers - force text format, assign, then restore the default style.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12"
$ws.Range("A2").Style = $ws.Range("B2").Style

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "3"
$ws.Range("H2").Style = $ws.Range("B2").Style

# --- Drop the old 3rd data row entirely (only one entry remains) ---
$ws.Rows.Item(3).Delete()
